$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting the existing rows 3-6 down to 4-7
$ws.Rows.Item(3).Insert()

# Fill the new row 3 with the new data record
$ws.Cells.Item(3, 1).Value = 3
$ws.Cells.Item(3, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(3, 3).Value = "Coquimbo"
$ws.Cells.Item(3, 4).Value = 44586
$ws.Cells.Item(3, 4).NumberFormat = $ws.Cells.Item(4, 4).NumberFormat
$ws.Cells.Item(3, 5).Value = 5
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100101
$ws.Cells.Item(3, 8).Value = "Berries"
$ws.Cells.Item(3, 9).Value = 100101004
$ws.Cells.Item(3, 10).Value = "Frambuesa"
$ws.Cells.Item(3, 11).Value = "Sin especificar"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 80
$ws.Cells.Item(3, 14).Value = 7000
$ws.Cells.Item(3, 15).Value = 7000
$ws.Cells.Item(3, 16).Value = 7000
$ws.Cells.Item(3, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(3, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(3, 19).Value = 3500
$ws.Cells.Item(3, 20).Value = 2
